# Adding more styles and a second sheet to the example xlsx

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after Sheet1; it becomes the active sheet,
# mirroring the commit's new "my_sheet 2" tab.
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "my_sheet 2"

# Row 1: "test this second" across A1:C1
$ws2.Range("A1").Value = "test"
$ws2.Range("B1").Value = "this"
$ws2.Range("C1").Value = "second"

# Row 2: "sheet" continues the sentence, and B2 demonstrates a wrapped
# cell containing an embedded line break.
$ws2.Range("A2").Value = "sheet"
$ws2.Range("B2").Value = "This one has a carriage`nReturn"
$ws2.Range("B2").WrapText = $true
$ws2.Rows.Item(2).RowHeight = 24

# Row 3: B3 demonstrates mixed run formatting - a strikethrough word in
# the middle of an otherwise plain sentence.
$ws2.Range("B3").Value = "this one is striked out"
$ws2.Range("B3").Characters(6, 14).Font.Strikethrough = $true

# Row 4: A4 demonstrates an underlined word.
$ws2.Range("A4").Value = "underlined"
$ws2.Range("A4").Font.Underline = $true

# Column B is widened so the wrapped text is readable.
$ws2.Columns.Item(2).ColumnWidth = 20.25

# Leave the selection where the author left it.
$ws2.Range("B6").Select() | Out-Null
